# Refresh the "Availability" workbook from its external source ("Wildnis u.
# Panorama") and then break the external link, matching what Excel does when
# you choose Data > Edit Links > Break Link after the linked values have
# changed: every formula that referenced the external workbook is replaced
# by its last-calculated (now static) value, and the external link
# definition is dropped from the workbook entirely.

$wb = $excel.ActiveWorkbook

# The refreshed values pulled from '[1]Wildnis u. Panorama'!$C$35:$AQ$35
# (displayed here as 24-minus-booked, i.e. what the formulas used to
# compute) -- identical sequence feeds both Sheet1!A5:AO5 and
# Availability!E2:E42.
$refreshed = @(24,24,21,23,24,22,24,24,22,23,14,24,24,16,22,17,24,22,22,24,24,24,22,17,20,22,24,20,22,24,24,13,22,24,24,24,22,24,22,24,24)

# --- Update the hidden "Sheet1" helper row (A5:AO5) ------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$row1 = New-Object 'object[,]' 1,41
for ($i = 0; $i -lt $refreshed.Count; $i++) {
    $row1[0, $i] = $refreshed[$i]
}
$ws1.Range("A5:AO5").Value = $row1

# --- Update the visible "Availability" sheet (E2:E42) ----------------------
$ws2 = $wb.Worksheets.Item("Availability")

for ($i = 0; $i -lt $refreshed.Count; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 5).Value = $refreshed[$i]
}

# Restore the default (no-override) selection on the Availability sheet.
$ws2.Activate() | Out-Null
$ws2.Range("A1").Select() | Out-Null

# --- Break the link to the external workbook --------------------------------
# Removes xl/externalLinks/externalLink1.xml + the workbook-level
# <externalReferences> entry, and (for any formula left referencing it)
# collapses the formula down to its cached value -- matching the rest of
# row 5 / column E which were already converted above.
foreach ($link in $wb.LinkSources()) {
    $wb.BreakLink($link, 1)
}
